{"js": "// The paragraph about Git's limitations ends with the phrase\n// \"...\u53ea\u80fd\u901a\u8fc7\u68c0\u67e5\u5927\u5c0f\u53d8\u5316\u533a\u522b\u6587\u6863\u5dee\u5f02\" where the final word \"\u5dee\u5f02\" needs to\n// become \"\u5dee\u522b\" (commit message: \"\u5dee \u6362\u6210 \u522b\"). There is an earlier, unrelated\n// occurrence of \"\u5dee\u5f02\" (\"\u5185\u5bb9\u5dee\u5f02\u53d8\u5316\") that must stay untouched, so we search\n// for enough trailing context to uniquely match only the final occurrence.\nconst body = context.document.body;\n\nconst target = \"\u53ea\u80fd\u901a\u8fc7\u68c0\u67e5\u5927\u5c0f\u53d8\u5316\u533a\u522b\u6587\u6863\u5dee\u5f02\";\nconst replacement = \"\u53ea\u80fd\u901a\u8fc7\u68c0\u67e5\u5927\u5c0f\u53d8\u5316\u533a\u522b\u6587\u6863\u5dee\u522b\";\n\nconst results = body.search(target, { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error('Could not find target text \"' + target + '\" in document body.');\n}\n\n// Replace only the first (and expected only) match in place, preserving the\n// run's existing formatting.\nresults.items[0].insertText(replacement, Word.InsertLocation.replace);\n\nawait context.sync();\n", "ps1": "# The paragraph about Git's limitations ends with the phrase\n# \"...\u53ea\u80fd\u901a\u8fc7\u68c0\u67e5\u5927\u5c0f\u53d8\u5316\u533a\u522b\u6587\u6863\u5dee\u5f02\" where the final word \"\u5dee\u5f02\" needs to\n# become \"\u5dee\u522b\" (commit message: \"\u5dee \u6362\u6210 \u522b\"). There is an earlier, unrelated\n# occurrence of \"\u5dee\u5f02\" (\"\u5185\u5bb9\u5dee\u5f02\u53d8\u5316\") that must stay untouched, so the find\n# text includes enough of the preceding sentence to uniquely match only the\n# final occurrence.\n$d = $word.ActiveDocument\n\n$searchRange = $d.Content\n$find = $searchRange.Find\n$find.ClearFormatting()\n$find.Text = \"\u53ea\u80fd\u901a\u8fc7\u68c0\u67e5\u5927\u5c0f\u53d8\u5316\u533a\u522b\u6587\u6863\u5dee\u5f02\"\n$find.Forward = $true\n$find.Wrap = 0  # wdFindStop - do not wrap around and risk a second match\n\n$find.Replacement.ClearFormatting()\n$find.Replacement.Text = \"\u53ea\u80fd\u901a\u8fc7\u68c0\u67e5\u5927\u5c0f\u53d8\u5316\u533a\u522b\u6587\u6863\u5dee\u522b\"\n\n# wdReplaceOne (1): replace just this single, uniquely-matched occurrence.\n$found = $find.Execute($find.Text, $false, $false, $false, $false, $false, $find.Forward, $find.Wrap, $false, $find.Replacement.Text, 1)\n\nif (-not $found) {\n  throw \"Could not find target text to replace in document.\"\n}\n"}
